$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 999.75
$ws.Range("I12").Value = 1109.7
$ws.Range("J12").Value = 450
$ws.Range("K12").Value = 1109.7
$ws.Range("L12").Value = 450
$ws.Range("M12").Value = -939.7
$ws.Range("N12").Value = -790

$ws.Range("H64").Value = 3378.6057
$ws.Range("I64").Value = 3016.0544
$ws.Range("J64").Value = 4624.875
$ws.Range("K64").Value = 3016.0544
$ws.Range("L64").Value = 4624.875
$ws.Range("M64").Value = -2768.0544
$ws.Range("N64").Value = -5120.875

$ws.Range("H67").Value = 3378.6057
$ws.Range("I67").Value = 3016.0544
$ws.Range("J67").Value = 4624.875
$ws.Range("K67").Value = 3016.0544
$ws.Range("L67").Value = 4624.875
$ws.Range("M67").Value = -2158.0544
$ws.Range("N67").Value = -6340.875

$ws.Range("H86").Value = 3415.889
$ws.Range("J86").Value = 3966.6667
$ws.Range("L86").Value = 3966.6667
$ws.Range("N86").Value = -6212.6667

$ws.Range("H89").Value = 3415.889
$ws.Range("J89").Value = 3966.6667
$ws.Range("L89").Value = 19833.3335
$ws.Range("N89").Value = -31065.3335

$ws.Range("H98").Value = 2458.4285
$ws.Range("I98").Value = 2534.8333
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2534.8333
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -1036.8333
$ws.Range("N98").Value = -4996

$ws.Range("H107").Value = 3972.2727
$ws.Range("J107").Value = 7264.4
$ws.Range("L107").Value = 7264.4
$ws.Range("N107").Value = -11104.4

$ws.Range("H122").Value = 2458.4285
$ws.Range("I122").Value = 2534.8333
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7604.499899999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -5154.499899999999
$ws.Range("N122").Value = -10900

$ws.Range("H131").Value = 12733.77
$ws.Range("I131").Value = 1294.9166
$ws.Range("K131").Value = 3884.7498
$ws.Range("M131").Value = 1155.2502

$ws.Range("H138").Value = 2819
$ws.Range("J138").Value = 3362.6216
$ws.Range("L138").Value = 10087.8648
$ws.Range("N138").Value = -20367.8648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 900.5
$ws.Range("I4").Value = 533.3333
$ws.Range("K4").Value = 533.3333
$ws.Range("M4").Value = -417.3333

$ws.Range("H5").Value = 152.11111
$ws.Range("I5").Value = 164.25
$ws.Range("J5").Value = 55
$ws.Range("K5").Value = 164.25
$ws.Range("L5").Value = 55
$ws.Range("M5").Value = -52.25
$ws.Range("N5").Value = -279

$ws.Range("H32").Value = 165079.86
$ws.Range("I32").Value = 218848.84
$ws.Range("J32").Value = 19587.295
$ws.Range("K32").Value = 218848.84
$ws.Range("L32").Value = 19587.295
$ws.Range("M32").Value = -218561.84
$ws.Range("N32").Value = -20161.295

$ws.Range("H61").Value = 8341176.5
$ws.Range("I61").Value = 9012
$ws.Range("J61").Value = 50002000
$ws.Range("K61").Value = 9012
$ws.Range("L61").Value = 50002000
$ws.Range("M61").Value = -8800
$ws.Range("N61").Value = -50002424

$ws.Range("H122").Value = 1332.5186
$ws.Range("I122").Value = 1177.1305
$ws.Range("K122").Value = 3531.3915
$ws.Range("M122").Value = -1081.3915

$ws.Range("H136").Value = 8341176.5
$ws.Range("I136").Value = 9012
$ws.Range("J136").Value = 50002000
$ws.Range("K136").Value = 27036
$ws.Range("L136").Value = 150006000
$ws.Range("M136").Value = -24486
$ws.Range("N136").Value = -150011100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 152.11111
$ws.Range("I4").Value = 164.25
$ws.Range("J4").Value = 55
$ws.Range("K4").Value = 164.25
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = -49.25
$ws.Range("N4").Value = -285

$ws.Range("H30").Value = 18471.5
$ws.Range("J30").Value = 18471.5
$ws.Range("L30").Value = 18471.5
$ws.Range("N30").Value = -18721.5

$ws.Range("H94").Value = 1365.1904
$ws.Range("I94").Value = 1414.9445
$ws.Range("J94").Value = 1066.6666
$ws.Range("K94").Value = 1414.9445
$ws.Range("L94").Value = 1066.6666
$ws.Range("M94").Value = -963.9445000000001
$ws.Range("N94").Value = -1968.6666

$ws.Range("H105").Value = 2456.465
$ws.Range("J105").Value = 3674.8333
$ws.Range("L105").Value = 3674.8333
$ws.Range("N105").Value = -7168.8333

$ws.Range("H107").Value = 957.7
$ws.Range("I107").Value = 771.2632
$ws.Range("K107").Value = 771.2632
$ws.Range("M107").Value = 1148.7368

$ws.Range("H134").Value = 5381691
$ws.Range("I134").Value = 5389.76
$ws.Range("J134").Value = 27782946
$ws.Range("K134").Value = 16169.28
$ws.Range("L134").Value = 83348838
$ws.Range("M134").Value = -13634.28
$ws.Range("N134").Value = -83353908

$ws.Range("H135").Value = 79666.664
$ws.Range("J135").Value = 79666.664
$ws.Range("L135").Value = 79666.664
$ws.Range("N135").Value = -89806.664

$ws.Range("H141").Value = 228889.67
$ws.Range("J141").Value = 228889.67
$ws.Range("L141").Value = 228889.67
$ws.Range("N141").Value = -239249.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1785.4286
$ws.Range("I105").Value = 1355.8
$ws.Range("K105").Value = 1355.8
$ws.Range("M105").Value = 391.2

$ws.Range("H107").Value = 267.2
$ws.Range("I107").Value = 267.2
$ws.Range("K107").Value = 267.2
$ws.Range("M107").Value = 1652.8

$ws.Range("H118").Value = 73483.375
$ws.Range("I118").Value = 65674
$ws.Range("J118").Value = 74599
$ws.Range("K118").Value = 65674
$ws.Range("L118").Value = 74599
$ws.Range("M118").Value = -64017
$ws.Range("N118").Value = -77913

$ws.Range("H132").Value = 2655.15
$ws.Range("I132").Value = 2595.077
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 7785.231000000001
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -5255.231000000001
$ws.Range("N132").Value = -20054

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 4049.8333
$ws.Range("I45").Value = 3100
$ws.Range("J45").Value = 4999.6665
$ws.Range("K45").Value = 9300
$ws.Range("L45").Value = 14998.9995
$ws.Range("M45").Value = -8768
$ws.Range("N45").Value = -16062.9995

$ws.Range("H121").Value = 2497.9167
$ws.Range("J121").Value = 2497.9167
$ws.Range("L121").Value = 7493.750100000001
$ws.Range("N121").Value = -10113.7501

$ws.Range("H122").Value = 1152978.9
$ws.Range("J122").Value = 1258
$ws.Range("L122").Value = 11322
$ws.Range("N122").Value = -16222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 38018.426
$ws.Range("I70").Value = 9541.947
$ws.Range("J70").Value = 76665.07000000001
$ws.Range("K70").Value = 9541.947
$ws.Range("L70").Value = 76665.07000000001
$ws.Range("M70").Value = -9271.947
$ws.Range("N70").Value = -77205.07000000001

$ws.Range("H73").Value = 38018.426
$ws.Range("I73").Value = 9541.947
$ws.Range("J73").Value = 76665.07000000001
$ws.Range("K73").Value = 9541.947
$ws.Range("L73").Value = 76665.07000000001
$ws.Range("M73").Value = -8605.947
$ws.Range("N73").Value = -78537.07000000001

$ws.Range("H102").Value = 2357.9
$ws.Range("I102").Value = 2286.5557
$ws.Range("K102").Value = 2286.5557
$ws.Range("M102").Value = -664.5556999999999

$ws.Range("H107").Value = 995.1
$ws.Range("I107").Value = 758.4091
$ws.Range("J107").Value = 1646
$ws.Range("K107").Value = 758.4091
$ws.Range("L107").Value = 1646
$ws.Range("M107").Value = 1161.5909
$ws.Range("N107").Value = -5486

$ws.Range("H122").Value = 46409.957
$ws.Range("I122").Value = 56497.26
$ws.Range("K122").Value = 169491.78
$ws.Range("M122").Value = -167041.78

$ws.Range("H132").Value = 4840.228
$ws.Range("I132").Value = 4312.163
$ws.Range("J132").Value = 8074.625
$ws.Range("K132").Value = 12936.489
$ws.Range("L132").Value = 24223.875
$ws.Range("M132").Value = -10406.489
$ws.Range("N132").Value = -29283.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4548.75
$ws.Range("I93").Value = 3731.8333
$ws.Range("K93").Value = 3731.8333
$ws.Range("M93").Value = -2483.8333

$ws.Range("H136").Value = 6759699
$ws.Range("I136").Value = 7815401
$ws.Range("J136").Value = 5955354.5
$ws.Range("K136").Value = 23446203
$ws.Range("L136").Value = 17866063.5
$ws.Range("M136").Value = -23443653
$ws.Range("N136").Value = -17871163.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13600
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 13600
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 13600
$ws.Range("N62").Value = -14848
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 13600
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 13600
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 68000
$ws.Range("N65").Value = -74240
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 3970745
$ws.Range("I132").Value = 4631898.5
$ws.Range("J132").Value = 3824.8333
$ws.Range("K132").Value = 13895695.5
$ws.Range("L132").Value = 11474.4999
$ws.Range("M132").Value = -13893165.5
$ws.Range("N132").Value = -16534.4999

$ws.Range("H136").Value = 14323260
$ws.Range("I136").Value = 2718464
$ws.Range("J136").Value = 200000000
$ws.Range("K136").Value = 8155392
$ws.Range("L136").Value = 600000000
$ws.Range("M136").Value = -8152842
$ws.Range("N136").Value = -600005100
